$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 157

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text
    if ($null -eq $val -or $val -eq "") { continue }

    $parts = $val -split ",\s*"
    $sysParts = @()
    $otherParts = @()
    foreach ($p in $parts) {
        if ($p.ToLower() -eq "system") {
            $sysParts += $p
        } else {
            $otherParts += $p
        }
    }

    if ($sysParts.Count -gt 0) {
        $newParts = $sysParts + $otherParts
        $newVal = $newParts -join ", "
        $cell.Value = $newVal
    }
}
